$d = $word.ActiveDocument

# The phrase "Thời gian cho mỗi câu hỏi là" occurs twice in the document
# (once under the WAR-mode parameter bullets, once under the "2 players
# answer together" bullets). The one we must change is the second
# occurrence: it is immediately preceded by the paragraph "Hai người chơi
# trả lời ... câu hỏi đó" and immediately followed by the bookmarked,
# yellow-highlighted paragraph "Mọi trợ giúp đều được áp dụng ở chế độ
# này". Anchor the search on that unique preceding sentence so we land on
# the right paragraph, then resolve the actual Paragraph object via the
# document's Paragraphs collection (Range.Paragraphs on a short / point
# range is not reliable in this engine).
$anchor = $d.Content
$null = $anchor.Find.Execute("Hai người chơi trả lời", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($anchor.End, $d.Content.End)
$null = $target.Find.Execute("Thời gian cho mỗi câu hỏi là", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hitPos = $target.Start

$paras = $d.Paragraphs
$targetIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Start -le $hitPos -and $p.Range.End -gt $hitPos) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    throw "Could not locate the 'Thời gian cho mỗi câu hỏi là' paragraph to edit"
}

$p = $d.Paragraphs.Item($targetIndex)
$paraStart = $p.Range.Start
$paraEnd = $p.Range.End

# $paraEnd points one past the paragraph mark; back off one character so
# the range below covers only the visible text -- "Thời gian cho mỗi câu
# hỏi là: " + tab + "20" + "s" -- without the paragraph mark itself.
$textRange = $d.Range($paraStart, $paraEnd - 1)
$textRange.Delete()

# Merge the now-empty bullet into the following paragraph by deleting its
# paragraph mark. As in Word, the surviving paragraph mark is the next
# paragraph's own mark (ListParagraph / numId 7 / yellow highlight), so
# the combined paragraph ends up with that formatting.
$emptyPara = $d.Paragraphs.Item($targetIndex)
$markRange = $d.Range($emptyPara.Range.Start, $emptyPara.Range.End)
$markRange.Delete()

# The _GoBack bookmark used to sit right before the paragraph mark we just
# removed; put it back at the same spot, which is now the very start of
# the merged paragraph (right before "Mọi trợ giúp...").
$bmRange = $d.Range($paraStart, $paraStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
